# "added bios and titles for November talk"
#
# The November 1st talks (row 26 = Madison Fitzpatrick / ghostwriting,
# row 27 = Andres Montealegre / psychology) were placeholders
# ("Title coming soon!" / " ") — fill in their real title + abstract now
# that the bios/titles are ready.
#
# Also clears the (duplicate, no-op) explicit style that had been applied
# to the "hashcode" column (I1:I23), returning those cells to the default
# style - matching the workbook's cleaned-up cellXfs table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Madison Fitzpatrick - ghostwriting talk
$ws.Range("G26").Value = "Ghostwriting: An Insider's View of an Invisible Profession"

# Andres Montealegre - psychology talk (abstract set before the ghostwriting
# abstract so the shared-string table ends up ordered the same way the
# authoring app produced it)
$ws.Range("H27").Value = " Ever wonder how others perceive you based on your actions? We all do. Shockingly, recent studies in psychology reveal we're often quite mistaken about these perceptions. For example, people think they're less liked by others than they actually are, or they assume others aren't interested in their unsolicited advice when in fact they are. Such findings could help us get along better with others. But how accurate are these studies? And is it possible that participants in these studies are just trying to appear 'nice' when answering research questions? In this talk, I'll share my own research on this topic and discuss if we should take these findings at face value. Join me to discuss the fascinating, messy world of how we're seen by others."

$ws.Range("H26").Value = " Ghostwriting often gets a bad rap as a cop-out for politicians and celebrities who are too lazy to write their own books. But the truth is, ghostwriting is more common than ever, and for many authors, it's no longer the shameful secret it once was. So, who actually hires ghostwriters? What do ghostwriters really do? And what's the future of this ancient profession in a world of generative AI? "

$ws.Range("G27").Value = "Can You Accurately Perceive How Others Think About You?"

# Clear the leftover duplicate style on the hashcode column
$ws.Range("I1:I23").ClearFormats()
